$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the header/body formatting from column M (which currently has the
#    full box border) onto the new column N, so N ends up with exactly the
#    style M used to have.
$ws.Range("M1:M11").Copy()
$ws.Range("N1:N11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Write the new "IP Address4" column values (this also creates the new
#    shared strings in the right order).
$ws.Range("N1").Value = "IP Address4"
$ws.Range("N2").Value = "44.201.90.191"
$ws.Range("N3").Value = "174.129.128.35"
$ws.Range("N4").Value = "3.83.45.68"
$ws.Range("N5").Value = "3.86.164.182"
$ws.Range("N6").Value = "34.230.74.252"
$ws.Range("N7").Value = "44.212.56.244"
$ws.Range("N8").Value = "3.89.31.66"
$ws.Range("N9").Value = "18.212.9.34"
$ws.Range("N10").Value = "44.206.229.188"
$ws.Range("N11").Value = "54.84.134.0"

# 3) Column N is now the rightmost bordered column, so remove the right edge
#    of the old rightmost column (M) to avoid a doubled border line.
$ws.Range("M1:M11").Borders.Item(10).LineStyle = -4142

# 4) Add one more (blank) row below the table with just the left/top/bottom
#    box border on N12 to close off the table visually.
$ws.Range("N12").Borders.Item(7).LineStyle = 1
$ws.Range("N12").Borders.Item(8).LineStyle = 1
$ws.Range("N12").Borders.Item(9).LineStyle = 1

# 5) Match the column width formatting applied to L:M onto the new N column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# 6) Restore the cursor/selection like the author left it.
$ws.Range("I22").Select()

Write-Output "done"
